# Player.xlsx - "Property" sheet updates
# 1) Queue-lock columns (E, "Save") for the MAXHP..MAXMP-style stat rows (44-67)
#    flipped from TRUE to FALSE.
# 2) GameID / GateID rows (76-77) lost their red/yellow "needs review" highlight
#    formatting now that the save-data layout changed; row 78 (GuildID) already
#    had the plain look these two rows are being brought in line with.
# 3) Selection cursor left on H78 (last edited cell) instead of E76.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# --- 2) Flip "Save" (column E) to FALSE for rows 44-67 ---
$ws.Range("E44:E67").Value = $false

# --- 3) Clear the red-font / yellow-fill highlight on rows 76 and 77 ---
$ws.Range("A76:J77").ClearFormats()

# A76 keeps a text number format (like the rows above it), the rest of the
# range reverts to the sheet's plain default formatting.
$ws.Range("A76").NumberFormat = "@"

# --- 1) Move the active selection to H78 ---
$null = $ws.Range("H78").Select()
